$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 156 (shifts existing rows 156-200 down to 157-201,
# and the used range grows from R200 to R201).
$ws.Rows.Item(156).Insert()

# Populate the newly inserted row with the new weekly price observation.
$ws.Range("A156").Value = 3
$ws.Range("B156").Value = "Femacal de La Calera"
$ws.Range("C156").Value = "Coquimbo"
$ws.Range("D156").Value = 44551
$ws.Range("E156").Value = 5
$ws.Range("F156").Value = 100112001
$ws.Range("G156").Value = "Berenjena"
$ws.Range("H156").Value = "Sin especificar"
$ws.Range("I156").Value = "Primera"
$ws.Range("J156").Value = 105
$ws.Range("K156").Value = 8000
$ws.Range("L156").Value = 9000
$ws.Range("M156").Value = 8524
$ws.Range("N156").Value = "$/caja 60 unidades"
$ws.Range("O156").Value = "Región de Arica y Parinacota"
$ws.Range("P156").Value = 142
$ws.Range("Q156").Value = 60
$ws.Range("R156").Value = "Hortaliza"
